{"js": "// Office.js edit script: \"Add some repository tests\"\n//\n// The commit rewrites two paragraphs of docs/TestDescription.docx:\n//   1. \"Testing The Service Layer\" paragraph:\n//        - \"such as \"Moq\" for xUnit)\" -> \"such as \"Moq\" for the xUnit library)\"\n//          (also drops the now-stale <w:proofErr/> spell-check bookmarks\n//          around \"Moq\"/\"xUnit\" since they get merged into plain runs)\n//        - \"I provided real implementations of\" -> \"I simply provided real implementations of\"\n//        - \"In this way, I verified the implementation of the data storage at\n//           the same time. Normally I would mock the interfaces\"\n//              -> \"Usually these would contain more complicated dependencies\n//                  and I would mock the interfaces\"\n//   2. \"Testing The Repository Layer\" paragraph:\n//        - \"Usually, this layer is implemented by ...\" -> \"Usually, the data\n//           storage is implemented by ...\"\n//        - \"I did not unit test the data storage separately since it is an\n//           implementation detail; and as stated before I would normally mock\n//           this layer or provide a fake.\"\n//              -> \"I added some simple tests for each repository method to\n//                  verify my data storage implementation.\"\n//\n// Rather than chaining many fragile range.search()/insertText() calls (Word's\n// search ranges can straddle existing runs in ways that leave orphaned\n// <w:proofErr/> markers behind), each paragraph is rewritten in one shot by\n// replacing its Range with the exact target WordprocessingML via\n// range.insertOoxml(..., Word.InsertLocation.replace). That mirrors exactly\n// what the authoritative diff shows for the run layout.\n\nconst W_NS = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\";\nconst W14_NS = \"http://schemas.microsoft.com/office/word/2010/wordml\";\n\n// Wrap a single <w:p>...</w:p> fragment in the pkg:package envelope that\n// Office.js's insertOoxml expects.\nfunction wrapParagraphOoxml(paragraphXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    `<w:document xmlns:w=\"${W_NS}\" xmlns:w14=\"${W14_NS}\">` +\n    `<w:body>${paragraphXml}</w:body>` +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nfunction run(text, attrs) {\n  const a = attrs ? \" \" + attrs : \"\";\n  // Word always needs xml:space=\"preserve\" whenever leading/trailing\n  // whitespace in the run text matters.\n  const preserve = /^\\s|\\s$/.test(text) ? ' xml:space=\"preserve\"' : \"\";\n  return `<w:r${a}><w:t${preserve}>${text}</w:t></w:r>`;\n}\n\n// --- Paragraph 1: \"Testing The Service Layer\" body paragraph -------------\nconst serviceLayerParagraphXml =\n  '<w:p w14:paraId=\"76D1B394\" w14:textId=\"4CF371E8\" w:rsidR=\"00172E92\" ' +\n  'w:rsidRDefault=\"00517EDF\" w:rsidP=\"00517EDF\">' +\n  run(\"I was not sure if I could use third party \") +\n  run(\"m\", 'w:rsidR=\"00172E92\"') +\n  run(\"ocking libraries\") +\n  run(\" (\", 'w:rsidR=\"00CA23A3\"') +\n  run(\"such as \\u201cMoq\\u201d for\") +\n  run(\" the\") +\n  run(\" xUnit\") +\n  run(\" library\") +\n  run(\") to mock the repository interfaces\", 'w:rsidR=\"00CA23A3\"') +\n  run(\". So \") +\n  run(\"instead,\", 'w:rsidR=\"00CA23A3\"') +\n  run(\" I \") +\n  run(\"simply \") +\n  run(\"provided real implementations of \") +\n  run(\"the repositories\", 'w:rsidR=\"00CA23A3\"') +\n  run(\" when testing the service layer. \") +\n  run(\"Usually\") +\n  run(\" \") +\n  run(\"these would contain more complicated dependencies and \") +\n  run(\"I would mock the interfaces\") +\n  run(\" to control the responses from them\", 'w:rsidR=\"0007770D\"') +\n  run(\". \", 'w:rsidR=\"00627842\"') +\n  \"</w:p>\";\n\n// --- Paragraph 2: \"Testing The Repository Layer\" body paragraph ----------\nconst repoLayerParagraphXml =\n  '<w:p w14:paraId=\"5642786F\" w14:textId=\"15536C78\" w:rsidR=\"00B67093\" ' +\n  'w:rsidRDefault=\"00172E92\" w:rsidP=\"00517EDF\">' +\n  run(\"Usually, \") +\n  run(\"the data storage\") +\n  run(\" is implemented by a third-party solution in which case the actual implementation shall already be \") +\n  run(\"unit \", 'w:rsidR=\"00BC38D5\"') +\n  run(\"tested. \") +\n  run(\"I added some simple tests for each repository method\") +\n  run(\" to verify my data storage implementation\") +\n  run(\".\") +\n  \"</w:p>\";\n\n// Locate the two target paragraphs by distinctive text they still contain\n// after the edit, rather than hard-coding paragraph indices.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet serviceLayerParagraph = null;\nlet repoLayerParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraphText = paragraphs.items[i].text;\n  if (paragraphText.indexOf(\"I was not sure if I could use third party\") !== -1) {\n    serviceLayerParagraph = paragraphs.items[i];\n  }\n  if (paragraphText.indexOf(\"this layer is implemented by a third-party solution\") !== -1) {\n    repoLayerParagraph = paragraphs.items[i];\n  }\n}\n\nif (!serviceLayerParagraph || !repoLayerParagraph) {\n  throw new Error(\"Could not locate the paragraphs to edit\");\n}\n\nserviceLayerParagraph\n  .getRange()\n  .insertOoxml(wrapParagraphOoxml(serviceLayerParagraphXml), Word.InsertLocation.replace);\nawait context.sync();\n\nrepoLayerParagraph\n  .getRange()\n  .insertOoxml(wrapParagraphOoxml(repoLayerParagraphXml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# PowerShell / Word COM edit script: \"Add some repository tests\"\n#\n# The commit rewrites two paragraphs of docs/TestDescription.docx:\n#   1. \"Testing The Service Layer\" paragraph:\n#        - 'such as \"Moq\" for xUnit)' -> 'such as \"Moq\" for the xUnit library)'\n#          (also drops the now-stale <w:proofErr/> spell-check bookmarks\n#          around \"Moq\"/\"xUnit\" since they get merged into plain runs)\n#        - \"I provided real implementations of\" -> \"I simply provided real implementations of\"\n#        - \"In this way, I verified the implementation of the data storage at\n#           the same time. Normally I would mock the interfaces\"\n#              -> \"Usually these would contain more complicated dependencies\n#                  and I would mock the interfaces\"\n#   2. \"Testing The Repository Layer\" paragraph:\n#        - \"Usually, this layer is implemented by ...\" -> \"Usually, the data\n#           storage is implemented by ...\"\n#        - \"I did not unit test the data storage separately since it is an\n#           implementation detail; and as stated before I would normally mock\n#           this layer or provide a fake.\"\n#              -> \"I added some simple tests for each repository method to\n#                  verify my data storage implementation.\"\n#\n# Rather than chaining many fragile Find/Replace calls (Find.Execute ranges\n# can straddle existing runs in ways that leave orphaned <w:proofErr/>\n# markers behind), each paragraph is rewritten in one shot by replacing its\n# Range with the exact target WordprocessingML via Range.InsertXML. That\n# mirrors exactly what the authoritative diff shows for the run layout.\n\nfunction New-Run {\n    param(\n        [string]$Text,\n        [string]$Attrs = \"\"\n    )\n    $openTag = \"<w:r>\"\n    if (-not [string]::IsNullOrEmpty($Attrs)) {\n        $openTag = \"<w:r $Attrs>\"\n    }\n    # Word always needs xml:space=\"preserve\" whenever leading/trailing\n    # whitespace in the run text matters.\n    $preserve = \"\"\n    if ($Text -match '^\\s' -or $Text -match '\\s$') {\n        $preserve = ' xml:space=\"preserve\"'\n    }\n    return \"$openTag<w:t$preserve>$Text</w:t></w:r>\"\n}\n\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n$w14Ns = 'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"'\n\n# --- Paragraph 1: \"Testing The Service Layer\" body paragraph -------------\n$serviceLayerRuns =\n    (New-Run \"I was not sure if I could use third party \") +\n    (New-Run \"m\" 'w:rsidR=\"00172E92\"') +\n    (New-Run \"ocking libraries\") +\n    (New-Run \" (\" 'w:rsidR=\"00CA23A3\"') +\n    (New-Run \"such as \u201cMoq\u201d for\") +\n    (New-Run \" the\") +\n    (New-Run \" xUnit\") +\n    (New-Run \" library\") +\n    (New-Run \") to mock the repository interfaces\" 'w:rsidR=\"00CA23A3\"') +\n    (New-Run \". So \") +\n    (New-Run \"instead,\" 'w:rsidR=\"00CA23A3\"') +\n    (New-Run \" I \") +\n    (New-Run \"simply \") +\n    (New-Run \"provided real implementations of \") +\n    (New-Run \"the repositories\" 'w:rsidR=\"00CA23A3\"') +\n    (New-Run \" when testing the service layer. \") +\n    (New-Run \"Usually\") +\n    (New-Run \" \") +\n    (New-Run \"these would contain more complicated dependencies and \") +\n    (New-Run \"I would mock the interfaces\") +\n    (New-Run \" to control the responses from them\" 'w:rsidR=\"0007770D\"') +\n    (New-Run \". \" 'w:rsidR=\"00627842\"')\n\n$serviceLayerXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId=\"76D1B394\" w14:textId=\"4CF371E8\" w:rsidR=\"00172E92\" w:rsidRDefault=\"00517EDF\" w:rsidP=\"00517EDF\">' + $serviceLayerRuns + '</w:p>'\n\n# --- Paragraph 2: \"Testing The Repository Layer\" body paragraph ----------\n$repoLayerRuns =\n    (New-Run \"Usually, \") +\n    (New-Run \"the data storage\") +\n    (New-Run \" is implemented by a third-party solution in which case the actual implementation shall already be \") +\n    (New-Run \"unit \" 'w:rsidR=\"00BC38D5\"') +\n    (New-Run \"tested. \") +\n    (New-Run \"I added some simple tests for each repository method\") +\n    (New-Run \" to verify my data storage implementation\") +\n    (New-Run \".\")\n\n$repoLayerXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId=\"5642786F\" w14:textId=\"15536C78\" w:rsidR=\"00B67093\" w:rsidRDefault=\"00172E92\" w:rsidP=\"00517EDF\">' + $repoLayerRuns + '</w:p>'\n\n# Locate the two target paragraphs by distinctive text they still contain\n# after the edit, rather than hard-coding paragraph indices.\n$d = $word.ActiveDocument\n\n$serviceLayerPara = $null\n$repoLayerPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*I was not sure if I could use third party*\") {\n        $serviceLayerPara = $p\n    }\n    if ($t -like \"*this layer is implemented by a third-party solution*\") {\n        $repoLayerPara = $p\n    }\n}\n\nif ($serviceLayerPara -eq $null -or $repoLayerPara -eq $null) {\n    throw \"Could not locate the paragraphs to edit\"\n}\n\n$serviceLayerPara.Range.InsertXML($serviceLayerXml)\n$repoLayerPara.Range.InsertXML($repoLayerXml)\n"}
